# Apply crypto price/volume updates from the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    @{Cell="D2"; Value="28.244.42"},
    @{Cell="E2"; Value="  +2.71%  "},
    @{Cell="D3"; Value="1.586.24"},
    @{Cell="E3"; Value="  +1.52%  "},
    @{Cell="E4"; Value="  +1.11%  "},
    @{Cell="D5"; Value="213.44"},
    @{Cell="E5"; Value="  +1.39%  "},
    @{Cell="E6"; Value="  +0.16%  "},
    @{Cell="E7"; Value="  +1.11%  "},
    @{Cell="D8"; Value="23.93"},
    @{Cell="E8"; Value="  +6.56%  "},
    @{Cell="E9"; Value="  -0.22%  "},
    @{Cell="D10"; Value="0.0598"},
    @{Cell="E10"; Value="  +0.19%  "},
    @{Cell="D11"; Value="0.0886"},
    @{Cell="E11"; Value="  +2.16%  "},
    @{Cell="D12"; Value="1.812.93"},
    @{Cell="E12"; Value="  +1.53%  "},
    @{Cell="D13"; Value="1.584.96"},
    @{Cell="E13"; Value="  +1.46%  "},
    @{Cell="D14"; Value="0.530"},
    @{Cell="E14"; Value="  +1.66%  "},
    @{Cell="E15"; Value="  -0.15%  "},
    @{Cell="D16"; Value="28.251.20"},
    @{Cell="E16"; Value="  +2.80%  "},
    @{Cell="D17"; Value="63.19"},
    @{Cell="E17"; Value="  +1.12%  "},
    @{Cell="D18"; Value="227.28"},
    @{Cell="E18"; Value="  +1.63%  "},
    @{Cell="E19"; Value="  -0.13%  "},
    @{Cell="D20"; Value="7.46"},
    @{Cell="E20"; Value="  -0.74%  "},
    @{Cell="E21"; Value="  +1.16%  "},
    @{Cell="E22"; Value="  -1.71%  "},
    @{Cell="E23"; Value="  -1.15%  "},
    @{Cell="E24"; Value="  +0.61%  "},
    @{Cell="D25"; Value="151.87"},
    @{Cell="E25"; Value="  +1.19%  "},
    @{Cell="E26"; Value="  -0.13%  "},
    @{Cell="E27"; Value="  -0.98%  "},
    @{Cell="E28"; Value="  -0.94%  "},
    @{Cell="E29"; Value="  +1.11%  "},
    @{Cell="E30"; Value="  -0.16%  "},
    @{Cell="E31"; Value="  +0.00%  "},
    @{Cell="E32"; Value="  -0.27%  "},
    @{Cell="E33"; Value="  -1.28%  "},
    @{Cell="D34"; Value="1.400.07"},
    @{Cell="E34"; Value="  -4.30%  "},
    @{Cell="D35"; Value="1.60"},
    @{Cell="E35"; Value="  -1.45%  "},
    @{Cell="E36"; Value="  -7.12%  "},
    @{Cell="E37"; Value="  +1.55%  "},
    @{Cell="D38"; Value="0.0167"},
    @{Cell="E38"; Value="  -0.44%  "},
    @{Cell="D39"; Value="2.54"},
    @{Cell="E39"; Value="  +8.86%  "},
    @{Cell="E40"; Value="  -0.30%  "},
    @{Cell="E41"; Value="  -0.48%  "},
    @{Cell="E42"; Value="  +1.14%  "},
    @{Cell="B43"; Value="RenderToken"},
    @{Cell="C43"; Value="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"},
    @{Cell="D43"; Value="1.87"},
    @{Cell="E43"; Value="  -1.24%  "},
    @{Cell="B44"; Value="FraxShare"},
    @{Cell="C44"; Value="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"},
    @{Cell="D44"; Value="5.59"},
    @{Cell="E44"; Value="  -1.94%  "},
    @{Cell="D45"; Value="0.979"},
    @{Cell="E45"; Value="  +0.36%  "},
    @{Cell="D46"; Value="64.20"},
    @{Cell="E46"; Value="  -1.37%  "},
    @{Cell="D47"; Value="1.722.83"},
    @{Cell="E47"; Value="  +1.27%  "},
    @{Cell="E48"; Value="  +2.12%  "},
    @{Cell="D49"; Value="86.95"},
    @{Cell="E49"; Value="  +0.48%  "},
    @{Cell="E50"; Value="  -1.16%  "},
    @{Cell="D51"; Value="0.0518"},
    @{Cell="E51"; Value="  -0.80%  "},
)

foreach ($edit in $edits) {
    $r = $ws.Range($edit.Cell)
    $r.NumberFormat = "@"
    $r.Value = $edit.Value
    $r.Style = "Normal"
}
